$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (A2:A5) to mirror the D column PK values / add REC row
$ws.Range("A2").Value = "AD.SEC.001.FON.02"

$ws.Range("A3").Value = "AD.SEC.001.FON.01"
$ws.Range("A3").NumberFormat = "@"

$ws.Range("A4").Value = "AD.SEC.001.FON.03"
$ws.Range("A5").Value = "RO.ACT"

# Add new row 6 data (REC line)
$ws.Range("A6").Value = "AD.SEC.014.FON.01"
$ws.Range("E6").Value = "RO.ACT.001.REC"

# Update the selection shown in the workbook to the new active range
$ws.Range("B12:B13").Select()
